# Ajout EDT S6 25-26.
# Fill in the room/location column (F) for each course/TD session that was
# previously left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "U3-110"
$ws.Range("F4").Value = "U3-110"
$ws.Range("F7").Value = "U3-4"
$ws.Range("F9").Value = "U3-109"
$ws.Range("F12").Value = "U3-Amphi"
$ws.Range("F13").Value = "U3-4"
$ws.Range("F16").Value = "U3-4"
$ws.Range("F19").Value = "U3-Amphi"
